$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32 (pushes the existing rows 32..51 down to 33..52,
# and carries the date-style formatting from the surrounding rows onto the new row)
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with its data
$ws.Cells.Item(32,1).Value = 11
$ws.Cells.Item(32,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(32,3).Value = "Bíobío"
$ws.Cells.Item(32,4).Value = 45062
$ws.Cells.Item(32,5).Value = 8
$ws.Cells.Item(32,6).Value = 100114007
$ws.Cells.Item(32,7).Value = "Jengibre"
$ws.Cells.Item(32,8).Value = "Sin especificar"
$ws.Cells.Item(32,9).Value = "Primera"
$ws.Cells.Item(32,10).Value = 30
$ws.Cells.Item(32,11).Value = 16000
$ws.Cells.Item(32,12).Value = 17000
$ws.Cells.Item(32,13).Value = 16333
$ws.Cells.Item(32,14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(32,15).Value = "Perú"
$ws.Cells.Item(32,16).Value = 1256
$ws.Cells.Item(32,17).Value = 13
$ws.Cells.Item(32,18).Value = "Hortaliza"
